$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G
$data = @(
    @(0.05294722152870939, 2.238670449449698, 21.89410068642309, 4.67911323718748, 4.788918458868004, 22),
    @(-0.06993565119730449, 2.484982260240081, 20.6686104960718, 4.54627435336582, 4.657993717879145, 21),
    @(-0.5679799476115557, 1.738661634217984, 9.770321410320426, 3.125751335330504, 3.153564430997775, 20),
    @(-0.1698634615297453, 1.407413280483496, 9.516255737529789, 3.084842903217243, 3.164566327120173, 19),
    @(-0.1701844293363445, 1.931438934375882, 13.53019059843001, 3.678340739848609, 3.780928157538495, 18),
    @(-0.2922294451516887, 1.859773159921934, 10.51533673933257, 3.242735995935002, 3.328935244021553, 17),
    @(-0.1996064344110716, 1.805775804098517, 11.52715139266361, 3.395165885882987, 3.500447005970189, 16),
    @(-0.08544148817198256, 1.883541101007153, 11.80030034152509, 3.435156523584491, 3.55462476826579, 15),
    @(-0.1451545775010554, 2.090306142858136, 13.68531860180392, 3.699367324530497, 3.836058499053441, 14),
    @(-0.1153228724130859, 2.085397217002122, 13.53761179153203, 3.679349370681184, 3.827706684449967, 13)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
